$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recommandations")

# Target data for rows 2-41 (A:G), reflecting the re-sorted/updated BRVM recommendations table.
$data = @(
    ,@('BRVM - SERVICES PUBLICS', 0, 5, 1746.93, 107.19, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - AUTRES SECTEURS', 0, 2, 1354.38, 669.2, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - DISTRIBUTION', 0, 2, 991.51, 494.78, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - TRANSPORT', 0, 2, 740.4400000000001, 370.83, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - AGRICULTURE', 0, 2, 673.1, 341.63, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - CONSOMMATION DISCRETIONNAIRE', 0, 3, 517.84, 170.29, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - INDUSTRIE  (**)', 0, 2, 512.28, 257.22, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - SERVICES FINANCIERS', 0, 3, 436.82, 145.91, '🟡 Observer', '➖ Neutre')
    ,@('BRVM-PRINCIPAL  (**)', 0, 2, 433.61, 217.65, '🟡 Observer', '➖ Neutre')
    ,@('BRVM-PRESTIGE', 0, 3, 429.89, 143.48, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - CONSOMMATION DE BASE  (**)', 0, 2, 428.95, 216.08, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - INDUSTRIELS', 0, 3, 396.85, 132.35, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - ENERGIE', 0, 3, 332.93, 110.99, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - FINANCES', 0, 2, 296, 148.41, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - TELECOMMUNICATIONS', 0, 3, 280.69, 94.05, '🟡 Observer', '➖ Neutre')
    ,@('BRVM-PRINCIPAL    (**)', 0, 1, 217.75, 217.75, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - CONSOMMATION DE BASE    (**)', 0, 1, 215.92, 215.92, '🟡 Observer', '➖ Neutre')
    ,@('BRVM – COMPOSITE TOTAL RETURN    (**)', 0, 1, 132.65, 132.65, '🟡 Observer', '➖ Neutre')
    ,@('EVIOSYS PACKAGING SIEM CI (SEMC)', 1, 0, 7.14, 7.14, '🟡 Observer', '➖ Neutre')
    ,@('SONATEL SN (SNTS)', 1, 0, 7.05, 7.05, '🟡 Observer', '➖ Neutre')
    ,@('ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)', 2, 1, 5.41, 6.33, '🟡 Observer', '👀 À surveiller')
    ,@('SAPH CI (SPHC)', 1, 0, 5.33, 5.33, '🟡 Observer', '➖ Neutre')
    ,@('BANK OF AFRICA BF (BOABF)', 1, 0, 5.07, 5.07, '🟡 Observer', '➖ Neutre')
    ,@('SOCIETE GENERALE COTE D''IVOIRE (SGBC)', 1, 0, 4.75, 4.75, '🟡 Observer', '➖ Neutre')
    ,@('CIE CI (CIEC)', 1, 0, 4.66, 4.66, '🟡 Observer', '➖ Neutre')
    ,@('ORAGROUP TOGO (ORGT)', 2, 1, 3.69, 3.54, '🟡 Observer', '👀 À surveiller')
    ,@('AFRICA GLOBAL LOGISTICS CI (SDSC)', 1, 0, 3.42, 3.42, '🟡 Observer', '➖ Neutre')
    ,@('NESTLE CI (NTLC)', 1, 0, 3.24, 3.24, '🟡 Observer', '➖ Neutre')
    ,@('SOCIETE IVOIRIENNE DE BANQUE  (SIBC)', 1, 0, 3.23, 3.23, '🟡 Observer', '➖ Neutre')
    ,@('SICOR CI (SICC)', 1, 1, 1, -5.71, '🟡 Observer', '👀 À surveiller')
    ,@('TOTAL', 0, 2, 0, 0, '🟡 Observer', '➖ Neutre')
    ,@('NEI-CEDA CI (NEIC)', 1, 1, -0.29, -7.5, '🟡 Observer', '👀 À surveiller')
    ,@('BERNABE CI (BNBC)', 0, 1, -3.19, -3.19, '🟡 Observer', '➖ Neutre')
    ,@('TRACTAFRIC MOTORS CI (PRSC)', 0, 1, -3.43, -3.43, '🟡 Observer', '➖ Neutre')
    ,@('BANK OF AFRICA ML (BOAM)', 0, 1, -3.5, -3.5, '🟡 Observer', '➖ Neutre')
    ,@('LOTERIE NATIONALE DU BENIN (LNBB)', 0, 1, -3.65, -3.65, '🟡 Observer', '➖ Neutre')
    ,@('CFAO MOTORS CI (CFAC)', 0, 2, -3.83, -2.45, '🟡 Observer', '➖ Neutre')
    ,@('SAFCA CI (SAFC)', 0, 1, -3.93, -3.93, '🟡 Observer', '➖ Neutre')
    ,@('SERVAIR ABIDJAN CI (ABJC)', 0, 1, -4.92, -4.92, '🟡 Observer', '➖ Neutre')
    ,@('SETAO CI (STAC)', 0, 3, -17.28, -3.4, '🔴 Vente', '⚠️ Risque de décrochage')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# Remove the now-obsolete trailing rows (table shrank from 44 to 41 rows; 3 tickers were dropped).
$ws.Rows.Item(44).Delete()
$ws.Rows.Item(43).Delete()
$ws.Rows.Item(42).Delete()

